# Auto-generated edit script: updates Tiamat_Profits market-price derived
# columns (H:N) across several leve rows on the ALC/ARM/BSM/CRP/CUL/LTW/WVR
# sheets, per the scheduled-runner price refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 408.33334
$ws.Range("I28").Value = 290
$ws.Range("J28").Value = 1000
$ws.Range("K28").Value = 290
$ws.Range("L28").Value = 1000
$ws.Range("M28").Value = 195
$ws.Range("N28").Value = -1970

$ws.Range("H33").Value = 324.375
$ws.Range("I33").Value = 116.388885
$ws.Range("J33").Value = 948.3333
$ws.Range("K33").Value = 116.388885
$ws.Range("L33").Value = 948.3333
$ws.Range("M33").Value = 112.611115
$ws.Range("N33").Value = -1406.3333

$ws.Range("H62").Value = 12184.193
$ws.Range("I62").Value = 13385
$ws.Range("J62").Value = 5940
$ws.Range("K62").Value = 13385
$ws.Range("L62").Value = 5940
$ws.Range("M62").Value = -12761
$ws.Range("N62").Value = -7188

$ws.Range("H65").Value = 12184.193
$ws.Range("I65").Value = 13385
$ws.Range("J65").Value = 5940
$ws.Range("K65").Value = 66925
$ws.Range("L65").Value = 29700
$ws.Range("M65").Value = -63805
$ws.Range("N65").Value = -35940

$ws.Range("H86").Value = 1564.5
$ws.Range("I86").Value = 1613.3077
$ws.Range("J86").Value = 930
$ws.Range("K86").Value = 1613.3077
$ws.Range("L86").Value = 930
$ws.Range("M86").Value = -490.3077000000001
$ws.Range("N86").Value = -3176

$ws.Range("H89").Value = 1564.5
$ws.Range("I89").Value = 1613.3077
$ws.Range("J89").Value = 930
$ws.Range("K89").Value = 8066.538500000001
$ws.Range("L89").Value = 4650
$ws.Range("M89").Value = -2450.538500000001
$ws.Range("N89").Value = -15882

$ws.Range("H92").Value = 650.17645
$ws.Range("I92").Value = 596.6429000000001
$ws.Range("J92").Value = 900
$ws.Range("K92").Value = 596.6429000000001
$ws.Range("L92").Value = 900
$ws.Range("M92").Value = 651.3570999999999
$ws.Range("N92").Value = -3396

$ws.Range("H98").Value = 4167271.2
$ws.Range("I98").Value = 581.4761999999999
$ws.Range("K98").Value = 581.4761999999999
$ws.Range("M98").Value = 916.5238000000001

$ws.Range("H100").Value = 4193.3447
$ws.Range("I100").Value = 2682.1365
$ws.Range("J100").Value = 8942.857
$ws.Range("K100").Value = 2682.1365
$ws.Range("L100").Value = 8942.857
$ws.Range("M100").Value = -2141.1365
$ws.Range("N100").Value = -10024.857

$ws.Range("H107").Value = 525.9
$ws.Range("I107").Value = 540
$ws.Range("K107").Value = 540
$ws.Range("M107").Value = 1380

$ws.Range("H122").Value = 4167271.2
$ws.Range("I122").Value = 581.4761999999999
$ws.Range("K122").Value = 1744.4286
$ws.Range("M122").Value = 705.5714000000003

$ws.Range("H138").Value = 1661.35
$ws.Range("J138").Value = 2127.3967
$ws.Range("L138").Value = 6382.1901
$ws.Range("N138").Value = -16662.1901

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 6835.625
$ws.Range("I2").Value = 705.0769
$ws.Range("J2").Value = 33401.332
$ws.Range("K2").Value = 705.0769
$ws.Range("L2").Value = 33401.332
$ws.Range("M2").Value = -592.0769
$ws.Range("N2").Value = -33627.332

$ws.Range("H102").Value = 2107.6667
$ws.Range("I102").Value = 2199.4
$ws.Range("K102").Value = 2199.4
$ws.Range("M102").Value = -577.4000000000001

$ws.Range("H116").Value = 6835.625
$ws.Range("I116").Value = 705.0769
$ws.Range("J116").Value = 33401.332
$ws.Range("K116").Value = 705.0769
$ws.Range("L116").Value = 33401.332
$ws.Range("M116").Value = 1588.9231
$ws.Range("N116").Value = -37989.332

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 6835.625
$ws.Range("I3").Value = 705.0769
$ws.Range("J3").Value = 33401.332
$ws.Range("K3").Value = 705.0769
$ws.Range("L3").Value = 33401.332
$ws.Range("M3").Value = -591.0769
$ws.Range("N3").Value = -33629.332

$ws.Range("H86").Value = 213481.22
$ws.Range("I86").Value = 1297
$ws.Range("J86").Value = 584803.5600000001
$ws.Range("K86").Value = 1297
$ws.Range("L86").Value = 584803.5600000001
$ws.Range("M86").Value = -174
$ws.Range("N86").Value = -587049.5600000001

$ws.Range("H89").Value = 213481.22
$ws.Range("I89").Value = 1297
$ws.Range("J89").Value = 584803.5600000001
$ws.Range("K89").Value = 6485
$ws.Range("L89").Value = 2924017.8
$ws.Range("M89").Value = -869
$ws.Range("N89").Value = -2935249.8

$ws.Range("H99").Value = 1411.1111
$ws.Range("I99").Value = 1257.1428
$ws.Range("J99").Value = 1950
$ws.Range("K99").Value = 1257.1428
$ws.Range("L99").Value = 1950
$ws.Range("M99").Value = 240.8571999999999
$ws.Range("N99").Value = -4946

$ws.Range("H105").Value = 1250
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 1250
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 1250
$ws.Range("M105").ClearContents()
$ws.Range("N105").Value = -4744

$ws.Range("H107").Value = 1125
$ws.Range("I107").Value = 1000
$ws.Range("J107").Value = 1166.6666
$ws.Range("K107").Value = 1000
$ws.Range("L107").Value = 1166.6666
$ws.Range("M107").Value = 920
$ws.Range("N107").Value = -5006.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1075
$ws.Range("I16").Value = 950
$ws.Range("J16").Value = 1200
$ws.Range("K16").Value = 950
$ws.Range("L16").Value = 1200
$ws.Range("M16").Value = -663
$ws.Range("N16").Value = -1774

$ws.Range("H105").Value = 1080.25
$ws.Range("I105").Value = 770
$ws.Range("J105").Value = 2011
$ws.Range("K105").Value = 770
$ws.Range("L105").Value = 2011
$ws.Range("M105").Value = 977
$ws.Range("N105").Value = -5505

$ws.Range("H113").Value = 1075
$ws.Range("I113").Value = 950
$ws.Range("J113").Value = 1200
$ws.Range("K113").Value = 950
$ws.Range("L113").Value = 1200
$ws.Range("M113").Value = 1220
$ws.Range("N113").Value = -5540

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 26259154
$ws.Range("I131").Value = 582.5
$ws.Range("J131").Value = 27524628
$ws.Range("K131").Value = 1747.5
$ws.Range("L131").Value = 82573884
$ws.Range("M131").Value = 3292.5
$ws.Range("N131").Value = -82583964

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1377.2941
$ws.Range("I22").Value = 756.75
$ws.Range("J22").Value = 1568.2307
$ws.Range("K22").Value = 756.75
$ws.Range("L22").Value = 1568.2307
$ws.Range("M22").Value = -461.75
$ws.Range("N22").Value = -2158.2307

$ws.Range("H27").Value = 1377.2941
$ws.Range("I27").Value = 756.75
$ws.Range("J27").Value = 1568.2307
$ws.Range("K27").Value = 756.75
$ws.Range("L27").Value = 1568.2307
$ws.Range("M27").Value = -649.75
$ws.Range("N27").Value = -1782.2307

$ws.Range("H40").Value = 38942.816
$ws.Range("I40").Value = 1707.2106
$ws.Range("J40").Value = 127377.375
$ws.Range("K40").Value = 1707.2106
$ws.Range("L40").Value = 127377.375
$ws.Range("M40").Value = -1571.2106
$ws.Range("N40").Value = -127649.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 189.92308
$ws.Range("I107").Value = 167
$ws.Range("J107").Value = 266.33334
$ws.Range("K107").Value = 501
$ws.Range("L107").Value = 799.0000200000001
$ws.Range("M107").Value = 1419
$ws.Range("N107").Value = -4639.00002

$ws.Range("H113").Value = 290.8125
$ws.Range("I113").Value = 297.3
$ws.Range("J113").Value = 280
$ws.Range("K113").Value = 891.9000000000001
$ws.Range("L113").Value = 840
$ws.Range("M113").Value = 1278.1
$ws.Range("N113").Value = -5180
